# Update filtered_output.xlsx data for both sheets (Neg_Change and Pos_Change)
# reflecting the latest filtered market snapshot.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# --- Neg_Change: grows from 9 data rows (2..9) to 11 data rows (2..12) ---
# Insert 3 additional rows before row 10 so rows 10-12 exist for the new data.
$ws1.Rows.Item(10).Insert()
$ws1.Rows.Item(10).Insert()
$ws1.Rows.Item(10).Insert()

$ws1.Cells.Item(2,1).Value = "INDUSINDBK"
$ws1.Cells.Item(2,2).Value = 773.45
$ws1.Cells.Item(2,3).Value = 779.8
$ws1.Cells.Item(2,4).Value = 765.7
$ws1.Cells.Item(2,5).Value = 769.75
$ws1.Cells.Item(2,6).Value = 1859573
$ws1.Cells.Item(2,7).Value = 4114205
$ws1.Cells.Item(2,8).Value = -0.5480115842550384
$ws1.Cells.Item(2,9).Value = "INDUSINDBK"
$ws1.Cells.Item(3,1).Value = "HAL"
$ws1.Cells.Item(3,2).Value = 4537
$ws1.Cells.Item(3,3).Value = 4568.5
$ws1.Cells.Item(3,4).Value = 4472.4
$ws1.Cells.Item(3,5).Value = 4543.9
$ws1.Cells.Item(3,6).Value = 1265929
$ws1.Cells.Item(3,7).Value = 2711914
$ws1.Cells.Item(3,8).Value = -0.5331972179058776
$ws1.Cells.Item(3,9).Value = "HAL"
$ws1.Cells.Item(4,1).Value = "LICI"
$ws1.Cells.Item(4,2).Value = 890
$ws1.Cells.Item(4,3).Value = 894.6
$ws1.Cells.Item(4,4).Value = 885
$ws1.Cells.Item(4,5).Value = 887
$ws1.Cells.Item(4,6).Value = 761942
$ws1.Cells.Item(4,7).Value = 1526696
$ws1.Cells.Item(4,8).Value = -0.5009209430037153
$ws1.Cells.Item(4,9).Value = "LICI"
$ws1.Cells.Item(5,1).Value = "JSL"
$ws1.Cells.Item(5,2).Value = 718.1
$ws1.Cells.Item(5,3).Value = 731
$ws1.Cells.Item(5,4).Value = 716.1
$ws1.Cells.Item(5,5).Value = 729.35
$ws1.Cells.Item(5,6).Value = 465289
$ws1.Cells.Item(5,7).Value = 945121
$ws1.Cells.Item(5,8).Value = -0.5076937238723931
$ws1.Cells.Item(5,9).Value = "JSL"
$ws1.Cells.Item(6,1).Value = "VOLTAS"
$ws1.Cells.Item(6,2).Value = 1270
$ws1.Cells.Item(6,3).Value = 1279.6
$ws1.Cells.Item(6,4).Value = 1258
$ws1.Cells.Item(6,5).Value = 1270
$ws1.Cells.Item(6,6).Value = 569544
$ws1.Cells.Item(6,7).Value = 1199406
$ws1.Cells.Item(6,8).Value = -0.5251449467486405
$ws1.Cells.Item(6,9).Value = "VOLTAS"
$ws1.Cells.Item(7,1).Value = "BSE"
$ws1.Cells.Item(7,2).Value = 2500
$ws1.Cells.Item(7,3).Value = 2514.4
$ws1.Cells.Item(7,4).Value = 2456
$ws1.Cells.Item(7,5).Value = 2479.8
$ws1.Cells.Item(7,6).Value = 2747105
$ws1.Cells.Item(7,7).Value = 5393999
$ws1.Cells.Item(7,8).Value = -0.4907108807398741
$ws1.Cells.Item(7,9).Value = "BSE"
$ws1.Cells.Item(8,1).Value = "POLICYBZR"
$ws1.Cells.Item(8,2).Value = 1852
$ws1.Cells.Item(8,3).Value = 1858.4
$ws1.Cells.Item(8,4).Value = 1825.1
$ws1.Cells.Item(8,5).Value = 1837
$ws1.Cells.Item(8,6).Value = 661217
$ws1.Cells.Item(8,7).Value = 1329799
$ws1.Cells.Item(8,8).Value = -0.5027692154979813
$ws1.Cells.Item(8,9).Value = "POLICYBZR"
$ws1.Cells.Item(9,1).Value = "BDL"
$ws1.Cells.Item(9,2).Value = 1591
$ws1.Cells.Item(9,3).Value = 1613.9
$ws1.Cells.Item(9,4).Value = 1572.1
$ws1.Cells.Item(9,5).Value = 1576
$ws1.Cells.Item(9,6).Value = 2039786
$ws1.Cells.Item(9,7).Value = 5047289
$ws1.Cells.Item(9,8).Value = -0.5958650277406347
$ws1.Cells.Item(9,9).Value = "BDL"
$ws1.Cells.Item(10,1).Value = "IGL"
$ws1.Cells.Item(10,2).Value = 206.31
$ws1.Cells.Item(10,3).Value = 207.33
$ws1.Cells.Item(10,4).Value = 203.35
$ws1.Cells.Item(10,5).Value = 204.3
$ws1.Cells.Item(10,6).Value = 858541
$ws1.Cells.Item(10,7).Value = 1747587
$ws1.Cells.Item(10,8).Value = -0.5087277486042183
$ws1.Cells.Item(10,9).Value = "IGL"
$ws1.Cells.Item(11,1).Value = "SJVN"
$ws1.Cells.Item(11,2).Value = 93.9
$ws1.Cells.Item(11,3).Value = 94
$ws1.Cells.Item(11,4).Value = 91.6
$ws1.Cells.Item(11,5).Value = 92
$ws1.Cells.Item(11,6).Value = 5003409
$ws1.Cells.Item(11,7).Value = 10371223
$ws1.Cells.Item(11,8).Value = -0.5175680823756272
$ws1.Cells.Item(11,9).Value = "SJVN"
$ws1.Cells.Item(12,1).Value = "PRESTIGE"
$ws1.Cells.Item(12,2).Value = 1630
$ws1.Cells.Item(12,3).Value = 1639.9
$ws1.Cells.Item(12,4).Value = 1596
$ws1.Cells.Item(12,5).Value = 1600
$ws1.Cells.Item(12,6).Value = 309115
$ws1.Cells.Item(12,7).Value = 668089
$ws1.Cells.Item(12,8).Value = -0.53731463921723
$ws1.Cells.Item(12,9).Value = "PRESTIGE"

# --- Pos_Change: shrinks from 19 data rows (2..20) to 16 data rows (2..17) ---
# Remove the trailing 3 rows (18-20) that are no longer part of the filtered set.
$ws2.Rows.Item(18).Delete()
$ws2.Rows.Item(18).Delete()
$ws2.Rows.Item(18).Delete()

$ws2.Cells.Item(2,1).Value = "SBIN"
$ws2.Cells.Item(2,2).Value = 822
$ws2.Cells.Item(2,3).Value = 828
$ws2.Cells.Item(2,4).Value = 819.1
$ws2.Cells.Item(2,5).Value = 825.5
$ws2.Cells.Item(2,6).Value = 7087573
$ws2.Cells.Item(2,7).Value = 4447508
$ws2.Cells.Item(2,8).Value = 0.5936054527614115
$ws2.Cells.Item(2,9).Value = "SBIN"
$ws2.Cells.Item(3,1).Value = "TCS"
$ws2.Cells.Item(3,2).Value = 3034
$ws2.Cells.Item(3,3).Value = 3055.9
$ws2.Cells.Item(3,4).Value = 3019.1
$ws2.Cells.Item(3,5).Value = 3022.5
$ws2.Cells.Item(3,6).Value = 3342174
$ws2.Cells.Item(3,7).Value = 2148166
$ws2.Cells.Item(3,8).Value = 0.555826691233359
$ws2.Cells.Item(3,9).Value = "TCS"
$ws2.Cells.Item(4,1).Value = "LT"
$ws2.Cells.Item(4,2).Value = 3700
$ws2.Cells.Item(4,3).Value = 3708.5
$ws2.Cells.Item(4,4).Value = 3641.7
$ws2.Cells.Item(4,5).Value = 3672.1
$ws2.Cells.Item(4,6).Value = 1769702
$ws2.Cells.Item(4,7).Value = 1227196
$ws2.Cells.Item(4,8).Value = 0.4420695634601156
$ws2.Cells.Item(4,9).Value = "LT"
$ws2.Cells.Item(5,1).Value = "HCLTECH"
$ws2.Cells.Item(5,2).Value = 1501
$ws2.Cells.Item(5,3).Value = 1516.2
$ws2.Cells.Item(5,4).Value = 1485.8
$ws2.Cells.Item(5,5).Value = 1489
$ws2.Cells.Item(5,6).Value = 4227538
$ws2.Cells.Item(5,7).Value = 2916128
$ws2.Cells.Item(5,8).Value = 0.4497093406050763
$ws2.Cells.Item(5,9).Value = "HCLTECH"
$ws2.Cells.Item(6,1).Value = "TATACONSUM"
$ws2.Cells.Item(6,2).Value = 1057.2
$ws2.Cells.Item(6,3).Value = 1059.4
$ws2.Cells.Item(6,4).Value = 1045.4
$ws2.Cells.Item(6,5).Value = 1045.4
$ws2.Cells.Item(6,6).Value = 767620
$ws2.Cells.Item(6,7).Value = 501792
$ws2.Cells.Item(6,8).Value = 0.5297573496588228
$ws2.Cells.Item(6,9).Value = "TATACONSUM"
$ws2.Cells.Item(7,1).Value = "PIDILITIND"
$ws2.Cells.Item(7,2).Value = 3061.1
$ws2.Cells.Item(7,3).Value = 3093.3
$ws2.Cells.Item(7,4).Value = 3055.1
$ws2.Cells.Item(7,5).Value = 3085.4
$ws2.Cells.Item(7,6).Value = 386023
$ws2.Cells.Item(7,7).Value = 255940
$ws2.Cells.Item(7,8).Value = 0.5082558412127842
$ws2.Cells.Item(7,9).Value = "PIDILITIND"
$ws2.Cells.Item(8,1).Value = "HAVELLS"
$ws2.Cells.Item(8,2).Value = 1467.4
$ws2.Cells.Item(8,3).Value = 1494.9
$ws2.Cells.Item(8,4).Value = 1463.1
$ws2.Cells.Item(8,5).Value = 1476
$ws2.Cells.Item(8,6).Value = 861774
$ws2.Cells.Item(8,7).Value = 557221
$ws2.Cells.Item(8,8).Value = 0.5465569316303586
$ws2.Cells.Item(8,9).Value = "HAVELLS"
$ws2.Cells.Item(9,1).Value = "LODHA"
$ws2.Cells.Item(9,2).Value = 1229.9
$ws2.Cells.Item(9,3).Value = 1239.9
$ws2.Cells.Item(9,4).Value = 1228.2
$ws2.Cells.Item(9,5).Value = 1233
$ws2.Cells.Item(9,6).Value = 1514599
$ws2.Cells.Item(9,7).Value = 964805
$ws2.Cells.Item(9,8).Value = 0.5698498660351056
$ws2.Cells.Item(9,9).Value = "LODHA"
$ws2.Cells.Item(10,1).Value = "INDIGO"
$ws2.Cells.Item(10,2).Value = 5978
$ws2.Cells.Item(10,3).Value = 6055
$ws2.Cells.Item(10,4).Value = 5955.5
$ws2.Cells.Item(10,5).Value = 5991.5
$ws2.Cells.Item(10,6).Value = 767646
$ws2.Cells.Item(10,7).Value = 510227
$ws2.Cells.Item(10,8).Value = 0.50451857702552
$ws2.Cells.Item(10,9).Value = "INDIGO"
$ws2.Cells.Item(11,1).Value = "GAIL"
$ws2.Cells.Item(11,2).Value = 173.06
$ws2.Cells.Item(11,3).Value = 174.84
$ws2.Cells.Item(11,4).Value = 172.87
$ws2.Cells.Item(11,5).Value = 173.52
$ws2.Cells.Item(11,6).Value = 8864011
$ws2.Cells.Item(11,7).Value = 6019021
$ws2.Cells.Item(11,8).Value = 0.4726665682010414
$ws2.Cells.Item(11,9).Value = "GAIL"
$ws2.Cells.Item(12,1).Value = "COFORGE"
$ws2.Cells.Item(12,2).Value = 1626
$ws2.Cells.Item(12,3).Value = 1649.9
$ws2.Cells.Item(12,4).Value = 1621.6
$ws2.Cells.Item(12,5).Value = 1642
$ws2.Cells.Item(12,6).Value = 1241786
$ws2.Cells.Item(12,7).Value = 873679
$ws2.Cells.Item(12,8).Value = 0.4213298019066499
$ws2.Cells.Item(12,9).Value = "COFORGE"
$ws2.Cells.Item(13,1).Value = "BHARATFORG"
$ws2.Cells.Item(13,2).Value = 1182.4
$ws2.Cells.Item(13,3).Value = 1186
$ws2.Cells.Item(13,4).Value = 1167.1
$ws2.Cells.Item(13,5).Value = 1183.2
$ws2.Cells.Item(13,6).Value = 638001
$ws2.Cells.Item(13,7).Value = 425145
$ws2.Cells.Item(13,8).Value = 0.5006668313163745
$ws2.Cells.Item(13,9).Value = "BHARATFORG"
$ws2.Cells.Item(14,1).Value = "KPITTECH"
$ws2.Cells.Item(14,2).Value = 1215.9
$ws2.Cells.Item(14,3).Value = 1230
$ws2.Cells.Item(14,4).Value = 1206
$ws2.Cells.Item(14,5).Value = 1206.5
$ws2.Cells.Item(14,6).Value = 704643
$ws2.Cells.Item(14,7).Value = 455057
$ws2.Cells.Item(14,8).Value = 0.5484719496678438
$ws2.Cells.Item(14,9).Value = "KPITTECH"
$ws2.Cells.Item(15,1).Value = "CONCOR"
$ws2.Cells.Item(15,2).Value = 541.55
$ws2.Cells.Item(15,3).Value = 542
$ws2.Cells.Item(15,4).Value = 531.1
$ws2.Cells.Item(15,5).Value = 532.5
$ws2.Cells.Item(15,6).Value = 1253153
$ws2.Cells.Item(15,7).Value = 894246
$ws2.Cells.Item(15,8).Value = 0.40135152966857
$ws2.Cells.Item(15,9).Value = "CONCOR"
$ws2.Cells.Item(16,1).Value = "HFCL"
$ws2.Cells.Item(16,2).Value = 71.95
$ws2.Cells.Item(16,3).Value = 72.46
$ws2.Cells.Item(16,4).Value = 70.35
$ws2.Cells.Item(16,5).Value = 71.1
$ws2.Cells.Item(16,6).Value = 7160882
$ws2.Cells.Item(16,7).Value = 5095065
$ws2.Cells.Item(16,8).Value = 0.4054544937110714
$ws2.Cells.Item(16,9).Value = "HFCL"
$ws2.Cells.Item(17,1).Value = "NBCC"
$ws2.Cells.Item(17,2).Value = 105.75
$ws2.Cells.Item(17,3).Value = 106.29
$ws2.Cells.Item(17,4).Value = 103
$ws2.Cells.Item(17,5).Value = 103.6
$ws2.Cells.Item(17,6).Value = 5824140
$ws2.Cells.Item(17,7).Value = 3763179
$ws2.Cells.Item(17,8).Value = 0.5476648865228042
$ws2.Cells.Item(17,9).Value = "NBCC"
